# "Apollo Bow and Sharur sprite"
#
# Applies, on the "List of items" worksheet:
#  - Row 29 (Sharur / Ninurta): mark the Status cell (E29) with the yellow
#    "in progress" highlight (sprite WIP marker), same as other weapon rows.
#  - Row 35 (Fragarach / Lugh): same yellow Status highlight.
#  - Row 54 (Neptune's trident): retype the item Name without the apostrophe.
#  - Row 59 (Apollo's Bow): mark Status yellow and set Category to "Weapon".
#  - Row 62 (Heracles' Club): fill in the missing Domain ("Hero").
#  - New row 63: Poseidon's trident entry (Greek / Poseidon / Water), with
#    the green "done" Status highlight and Category "Weapon".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List of items")

# Colors used elsewhere in the sheet for the Status column (E):
#   yellow = 65535   (work-in-progress)
#   green  = 5287936 (done / sprite exists)

# Row 29 - Sharur: add yellow status highlight
$ws.Cells.Item(29, 5).Interior.Color = 65535

# Row 35 - Fragarach: add yellow status highlight
$ws.Cells.Item(35, 5).Interior.Color = 65535

# Row 54 - Neptune's trident -> Neptunes trident
$ws.Cells.Item(54, 1).Value = "Neptunes trident"

# Row 59 - Apollo's Bow: add yellow status highlight + Category "Weapon"
$ws.Cells.Item(59, 5).Interior.Color = 65535
$ws.Cells.Item(59, 6).Value = "Weapon"

# Row 62 - Heracles' Club: fill in Domain
$ws.Cells.Item(62, 4).Value = "Hero"

# Row 63 (new) - Poseidons trident
$ws.Cells.Item(63, 1).Value = "Poseidons trident"
$ws.Cells.Item(63, 2).Value = "Greek"
$ws.Cells.Item(63, 3).Value = "Poseidon"
$ws.Cells.Item(63, 4).Value = "Water"
$ws.Cells.Item(63, 5).Interior.Color = 5287936
$ws.Cells.Item(63, 6).Value = "Weapon"

# Restore view state: scroll so row 22 is at the top and select C69
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 22 } catch {}
$ws.Range("C69").Select()
